$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.232.26"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.568.70"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.Value = "'584.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$c = $ws.Range("D6")
$c.Value = "'149.18"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  +3.80%  "
$c = $ws.Range("D10")
$c.Value = "'5.60"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.31%  "
$c = $ws.Range("D13")
$c.Value = "'27.99"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "3.028.62"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "63.121.03"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "2.563.03"
$ws.Range("E17").Value = "  +0.27%  "
$c = $ws.Range("D18")
$c.Value = "'11.47"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$c = $ws.Range("D19")
$c.Value = "'341.68"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "
$c = $ws.Range("D20")
$c.Value = "'4.41"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "
$c = $ws.Range("D21")
$c.Value = "'6.87"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("E22").Value = "  -0.31%  "
$c = $ws.Range("D23")
$c.Value = "'66.21"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'1.66"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.95%  "
$ws.Range("D25").Value = "2.688.63"
$ws.Range("E25").Value = "  +0.85%  "
$c = $ws.Range("D26")
$c.Value = "'0.171"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$c = $ws.Range("D27")
$c.Value = "'8.32"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +15.36%  "
$c = $ws.Range("D28")
$c.Value = "'8.59"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("D32").Value = "0.0₃0833"
$ws.Range("E32").Value = "  +1.76%  "
$c = $ws.Range("D33")
$c.Value = "'177.64"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$c = $ws.Range("D34")
$c.Value = "'441.91"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.33%  "
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  +1.99%  "
$ws.Range("E37").Value = "  +2.39%  "
$c = $ws.Range("D40")
$c.Value = "'1.77"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$c = $ws.Range("D41")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "
$c = $ws.Range("D42")
$c.Value = "'152.49"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "
$c = $ws.Range("D43")
$c.Value = "'3.83"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.32%  "
$c = $ws.Range("D44")
$c.Value = "'21.57"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.19%  "
$c = $ws.Range("D45")
$c.Value = "'0.0557"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.47%  "
$c = $ws.Range("D46")
$c.Value = "'0.610"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "
$c = $ws.Range("D47")
$c.Value = "'0.0979"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("E48").Value = "  +2.99%  "
$c = $ws.Range("D49")
$c.Value = "'18.55"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.25%  "
